$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-01 Saturday", 2) | Out-Null
$d.Content.Find.Execute("173÷4=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "149÷7=21, 2", 2) | Out-Null
$d.Content.Find.Execute("151÷2=75, 1", $true, $false, $false, $false, $false, $true, 1, $false, "714÷6=119, 0", 2) | Out-Null
$d.Content.Find.Execute("524÷6=87, 2", $true, $false, $false, $false, $false, $true, 1, $false, "517÷7=73, 6", 2) | Out-Null
$d.Content.Find.Execute("635÷9=70, 5", $true, $false, $false, $false, $false, $true, 1, $false, "288÷2=144, 0", 2) | Out-Null
$d.Content.Find.Execute("902÷2=451, 0", $true, $false, $false, $false, $false, $true, 1, $false, "455÷9=50, 5", 2) | Out-Null
$d.Content.Find.Execute("150÷3=50, 0", $true, $false, $false, $false, $false, $true, 1, $false, "328÷7=46, 6", 2) | Out-Null
$d.Content.Find.Execute("810÷3=270, 0", $true, $false, $false, $false, $false, $true, 1, $false, "325÷9=36, 1", 2) | Out-Null
$d.Content.Find.Execute("889÷2=444, 1", $true, $false, $false, $false, $false, $true, 1, $false, "480÷7=68, 4", 2) | Out-Null
$d.Content.Find.Execute("139÷3=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "397÷9=44, 1", 2) | Out-Null
$d.Content.Find.Execute("690÷5=138, 0", $true, $false, $false, $false, $false, $true, 1, $false, "724÷4=181, 0", 2) | Out-Null
$d.Content.Find.Execute("456÷5=91, 1", $true, $false, $false, $false, $false, $true, 1, $false, "889÷4=222, 1", 2) | Out-Null
$d.Content.Find.Execute("466÷3=155, 1", $true, $false, $false, $false, $false, $true, 1, $false, "739÷9=82, 1", 2) | Out-Null
$d.Content.Find.Execute("342÷7=48, 6", $true, $false, $false, $false, $false, $true, 1, $false, "861÷6=143, 3", 2) | Out-Null
$d.Content.Find.Execute("132÷9=14, 6", $true, $false, $false, $false, $false, $true, 1, $false, "106÷2=53, 0", 2) | Out-Null
$d.Content.Find.Execute("346÷7=49, 3", $true, $false, $false, $false, $false, $true, 1, $false, "815÷7=116, 3", 2) | Out-Null
$d.Content.Find.Execute("140÷3=46, 2", $true, $false, $false, $false, $false, $true, 1, $false, "963÷3=321, 0", 2) | Out-Null
$d.Content.Find.Execute("401÷2=200, 1", $true, $false, $false, $false, $false, $true, 1, $false, "821÷6=136, 5", 2) | Out-Null
$d.Content.Find.Execute("137÷9=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "896÷9=99, 5", 2) | Out-Null
$d.Content.Find.Execute("916÷5=183, 1", $true, $false, $false, $false, $false, $true, 1, $false, "798÷5=159, 3", 2) | Out-Null
$d.Content.Find.Execute("912÷8=114, 0", $true, $false, $false, $false, $false, $true, 1, $false, "559÷2=279, 1", 2) | Out-Null
$d.Content.Find.Execute("745÷2=372, 1", $true, $false, $false, $false, $false, $true, 1, $false, "301÷9=33, 4", 2) | Out-Null
$d.Content.Find.Execute("260÷4=65, 0", $true, $false, $false, $false, $false, $true, 1, $false, "217÷7=31, 0", 2) | Out-Null
$d.Content.Find.Execute("843÷9=93, 6", $true, $false, $false, $false, $false, $true, 1, $false, "584÷9=64, 8", 2) | Out-Null
$d.Content.Find.Execute("702÷5=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "372÷2=186, 0", 2) | Out-Null
$d.Content.Find.Execute("269÷7=38, 3", $true, $false, $false, $false, $false, $true, 1, $false, "229÷6=38, 1", 2) | Out-Null
